$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new firm's data (Lawyers column is stored as text,
# so prefix the numeric-looking value with an apostrophe to keep it text)
$ws.Range("A2").Value = "Njord Law"
$ws.Range("B2").Value = "9sec"
$ws.Range("C2").Value = "'1"

# Clear row 3 (firm slot emptied) - use an apostrophe so the cells stay
# text/empty-string (matching the other blank rows) instead of becoming
# numeric blank cells
$ws.Range("A3").Value = "'"
$ws.Range("B3").Value = "'"
$ws.Range("C3").Value = "'"
